$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New stock movement rows (auto fabric consumption processing) appended
# to the bottom of the existing stock_movement sheet.
$newRows = @(
    @(46002, "F001", "이태리 순모 네이비", "OUT", 2.7, "m", "2025-3811-01",  "자동 원단 소요 처리", -2.7),
    @(46020, "F001", "이태리 순모 네이비", "OUT", 2.7, "m", "2025-10000-01", "자동 원단 소요 처리", -2.7),
    @(46020, "F001", "이태리 순모 네이비", "OUT", 2.7, "m", "2025-10000-01", "자동 원단 소요 처리", -2.7),
    @(46020, "F001", "이태리 순모 네이비", "OUT", 2.7, "m", "2025-10000-02", "자동 원단 소요 처리", -2.7)
)

$startRow = 8
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
}
